$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.173.73"
$ws.Range("E2").Value = "  -4.91%  "
$ws.Range("D3").Value = "3.257.54"
$ws.Range("E3").Value = "  -7.45%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'596.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.99%  "
$ws.Range("D6").Value = "'150.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.72%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.248.66"
$ws.Range("E9").Value = "  -11.55%  "
$ws.Range("E10").Value = "  -14.03%  "
$ws.Range("D11").Value = "'6.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.16%  "
$ws.Range("D12").Value = "'0.505"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -14.00%  "
$ws.Range("D13").Value = "'38.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -17.74%  "
$ws.Range("E14").Value = "  -12.63%  "
$ws.Range("D15").Value = "3.777.84"
$ws.Range("E15").Value = "  -7.67%  "
$ws.Range("D16").Value = "67.224.58"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").Value = "3.257.31"
$ws.Range("E17").Value = "  -7.37%  "
$ws.Range("D18").Value = "'536.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -11.95%  "
$ws.Range("E19").Value = "  -6.29%  "
$ws.Range("D20").Value = "'7.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -14.06%  "
$ws.Range("D21").Value = "'15.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -14.88%  "
$ws.Range("D22").Value = "'0.761"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -13.85%  "
$ws.Range("D23").Value = "'7.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -13.59%  "
$ws.Range("D24").Value = "'85.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -12.39%  "
$ws.Range("D25").Value = "'13.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -13.32%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'3.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -12.78%  "
$ws.Range("D28").Value = "'29.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -12.69%  "
$ws.Range("D29").Value = "'8.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.70%  "
$ws.Range("D30").Value = "'2.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -17.26%  "
$ws.Range("D31").Value = "'2.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.67%  "
$ws.Range("E32").Value = "  -12.59%  "
$ws.Range("E33").Value = "  -17.78%  "
$ws.Range("D34").Value = "'540.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -14.84%  "
$ws.Range("D35").Value = "'5.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -16.62%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "'0.0450"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.23%  "
$ws.Range("D38").Value = "'53.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.40%  "
$ws.Range("D39").Value = "'0.0853"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -13.91%  "
$ws.Range("D40").Value = "'9.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -15.82%  "
$ws.Range("D41").Value = "'0.128"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.11%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -20.81%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.918.19"
$ws.Range("E43").Value = "  -12.83%  "
$ws.Range("D44").Value = "'0.262"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -15.95%  "
$ws.Range("E45").Value = "  -19.51%  "
$ws.Range("D46").Value = "'26.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -16.93%  "
$ws.Range("D47").Value = "'2.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -14.77%  "
$ws.Range("D49").Value = "'127.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.38%  "
$ws.Range("E51").Value = "  -12.75%  "
